# edit.ps1
#
# Applies the two substantive content changes captured by the commit:
#
#  1. The "datetimeFigureOut" auto-date fields (slide master, all slide
#     layouts, handout master, notes master) show a cached value of
#     "7/8/2021" and need to show "7/14/2021" instead.
#
#  2. On the notes pages that document the find_me_model.drawio diagram,
#     the two hyperlinked runs
#         "find_me_model.drawio" + " - diagrams.net"
#     (both pointing at the same diagrams.net hyperlink) collapse into a
#     single, non-hyperlinked run containing the bare URL:
#         "https://app.diagrams.net/#G17DjBDN78j0e-FsorBzU61J3DlUEmU9Wt"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update every cached "7/8/2021" date placeholder to "7/14/2021".
# ---------------------------------------------------------------------

function Update-DateShapes($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "7/8/2021") {
                $sh.TextFrame2.TextRange.Text = "7/14/2021"
            }
        }
    }
}

Update-DateShapes $p.SlideMaster

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    Update-DateShapes $p.SlideMaster.CustomLayouts.Item($li)
}

Update-DateShapes $p.HandoutMaster
Update-DateShapes $p.NotesMaster

# ---------------------------------------------------------------------
# 2) Replace the hyperlinked "find_me_model.drawio - diagrams.net" notes
#    text with the plain diagrams.net URL, wherever it appears.
# ---------------------------------------------------------------------

$targetUrl = "https://app.diagrams.net/#G17DjBDN78j0e-FsorBzU61J3DlUEmU9Wt"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    $notes = $null
    try {
        $notes = $s.NotesPage
    } catch {
        $notes = $null
    }

    if ($notes -ne $null) {
        for ($shi = 1; $shi -le $notes.Shapes.Count; $shi++) {
            $sh = $notes.Shapes.Item($shi)
            if ($sh.HasTextFrame) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t.Contains("find_me_model.drawio")) {
                    $sh.TextFrame.TextRange.Text = $targetUrl
                }
            }
        }
    }
}
